$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price values so they keep their original
# string representation (leading/trailing zeros, fixed decimal places) instead of
# being auto-converted to actual numbers by Excel.
$textCells = @("D4","D5","D8","D9","D11","D12","D13","D14","D15","D18","D20","D21","D23","D25","D26","D27","D28","D30","D31","D33","D34","D35","D36","D38","D39","D40","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.379.05"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "1.849.17"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "240.15"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.07630"
$ws.Range("D9").Value = "0.2904"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  +0.97%  "
$ws.Range("D11").Value = "0.07739"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "5.032"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "0.6783"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "0.00001055"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Value = "83.18"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "29.413.87"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "227.45"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("D21").Value = "7.504"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "158.60"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").Value = "8.402"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").Value = "17.68"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "1.379"
$ws.Range("E27").Value = "  +5.96%  "
$ws.Range("D28").Value = "1.459"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "4.108"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").Value = "4.067"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("D33").Value = "1.162"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").Value = "0.6986"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").Value = "2.579"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "0.01805"
$ws.Range("E36").Value = "  +0.43%  "
$ws.Range("D37").Value = "1.230.64"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "2.712"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "6.373"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "0.9011"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "101.50"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "66.05"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "7.202"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("D46").Value = "0.4009"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Value = "9.028"
$ws.Range("E47").Value = "  -0.37%  "
$ws.Range("D48").Value = "1.679"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "0.1133"
$ws.Range("E49").Value = "  +0.97%  "
$ws.Range("D50").Value = "0.05700"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").Value = "0.4628"
$ws.Range("E51").Value = "  +0.07%  "

Write-Host "Updated cryptos list"
